$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '28.446.98'; E = '  +1.53%  ' },
    @{ Row = 3; D = '1.885.86'; E = '  +0.88%  ' },
    @{ Row = 4; D = '1.019'; E = '  +1.62%  ' },
    @{ Row = 5; D = '317.00'; E = '  +1.54%  ' },
    @{ Row = 6; D = '1.019'; E = '  +1.77%  ' },
    @{ Row = 7; D = '0.5146'; E = '  +0.05%  ' },
    @{ Row = 8; D = '0.3911'; E = '  +1.64%  ' },
    @{ Row = 9; D = '0.08272'; E = '  -0.28%  ' },
    @{ Row = 10; D = '1.123'; E = '  +1.01%  ' },
    @{ Row = 11; B = 'OKB'; C = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; D = '42.17'; E = '  +1.67%  ' },
    @{ Row = 12; B = 'Polkadot'; C = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D = '6.242'; E = '  +0.58%  ' },
    @{ Row = 13; B = 'WrappedEther'; C = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D = '1.890.22'; E = '  +0.68%  ' },
    @{ Row = 14; B = 'Solana'; C = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D = '20.27'; E = '  -1.43%  ' },
    @{ Row = 15; B = 'Chainlink'; C = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D = '7.241'; E = '  -0.62%  ' },
    @{ Row = 16; B = 'BinanceUSD'; C = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D = '1.022'; E = '  +1.92%  ' },
    @{ Row = 17; B = 'ShibaInu'; C = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D = '0.00001103'; E = '  +0.40%  ' },
    @{ Row = 18; B = 'Litecoin'; C = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D = '90.96'; E = '  +0.24%  ' },
    @{ Row = 19; B = 'TRON'; C = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D = '0.06752'; E = '  +1.58%  ' },
    @{ Row = 20; B = 'Avalanche'; C = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D = '17.76'; E = '  +0.22%  ' },
    @{ Row = 21; B = 'Dai'; C = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'; D = '1.020'; E = '  +1.82%  ' },
    @{ Row = 22; B = 'Uniswap'; C = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D = '6.025'; E = '  -0.05%  ' },
    @{ Row = 23; B = 'WrappedBTC'; C = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D = '28.484.28'; E = '  +1.52%  ' },
    @{ Row = 24; B = 'Cosmos'; C = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D = '11.12'; E = '  +0.22%  ' },
    @{ Row = 25; B = 'Toncoin'; C = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D = '2.252'; E = '  +0.18%  ' },
    @{ Row = 26; B = 'WrappedliquidstakedEther2.0'; C = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; D = '2.110.23'; E = '  +1.74%  ' },
    @{ Row = 27; B = 'Monero'; C = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D = '161.18'; E = '  +2.39%  ' },
    @{ Row = 28; B = 'EthereumClassic'; C = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D = '20.70'; E = '  +0.77%  ' },
    @{ Row = 29; B = 'LidoDAOToken'; C = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D = '2.425'; E = '  -3.66%  ' },
    @{ Row = 30; B = 'BitcoinCash'; C = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D = '125.54'; E = '  +0.46%  ' },
    @{ Row = 31; B = 'Stellar'; C = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D = '0.1058'; E = '  -0.68%  ' },
    @{ Row = 32; B = 'ImmutableX'; C = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D = '1.040'; E = '  +0.76%  ' },
    @{ Row = 33; B = 'Filecoin'; C = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D = '5.897'; E = '  +1.26%  ' },
    @{ Row = 34; B = 'HuobiToken'; C = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; D = '3.651'; E = '  +1.79%  ' },
    @{ Row = 35; B = 'FraxShare'; C = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D = '9.332'; E = '  -1.53%  ' },
    @{ Row = 36; B = 'VeChain'; C = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D = '0.02437'; E = '  +0.86%  ' },
    @{ Row = 37; B = 'Hedera'; C = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; D = '0.06573'; E = '  +0.69%  ' },
    @{ Row = 38; B = 'Algorand'; C = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; D = '0.2192'; E = '  -0.21%  ' },
    @{ Row = 39; B = 'TrustWalletToken'; C = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D = '1.262'; E = '  +4.22%  ' },
    @{ Row = 40; B = 'TheSandbox'; C = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; D = '0.6501'; E = '  -0.78%  ' },
    @{ Row = 41; B = 'ARBITRUM'; C = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D = '1.184'; E = '  -1.76%  ' },
    @{ Row = 42; B = 'InternetComputer(DFINITY)'; C = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D = '4.972'; E = '  -1.04%  ' },
    @{ Row = 43; B = 'Aptos'; C = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D = '11.16'; E = '  -0.51%  ' },
    @{ Row = 44; B = 'Decentraland'; C = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'; D = '0.6097'; E = '  -0.55%  ' },
    @{ Row = 45; B = 'EnergySwap'; C = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D = '13.21'; E = '  +0.64%  ' },
    @{ Row = 46; B = 'PancakeSwap'; C = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; D = '3.722'; E = '  +1.24%  ' },
    @{ Row = 47; B = 'WEMIXTOKEN'; C = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; D = '1.290'; E = '  +0.93%  ' },
    @{ Row = 48; B = 'NEARProtocol'; C = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D = '2.004'; E = '  -0.92%  ' },
    @{ Row = 49; B = 'EOS'; C = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'; D = '1.218'; E = '  +0.24%  ' },
    @{ Row = 50; B = 'Quant'; C = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; D = '122.49'; E = '  +1.32%  ' },
    @{ Row = 51; B = 'Cronos'; C = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; D = '0.06923'; E = '  +1.54%  ' }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("B")) { $ws.Range("B" + $r).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Range("C" + $r).Value = $u.C }
    if ($u.ContainsKey("D")) {
        $ws.Range("D" + $r).NumberFormat = "@"
        $ws.Range("D" + $r).Value = $u.D
        $ws.Range("D" + $r).Style = "Normal"
    }
    if ($u.ContainsKey("E")) { $ws.Range("E" + $r).Value = $u.E }
}